$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not auto-converted to a number/date
# by COM's locale-aware parsing) while leaving the cell's style index
# unchanged from the sheet's default (matches the original inlineStr cells,
# which carry no explicit style).
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 6: additional-shares count changes from 3,766 to 5,022
Set-TextValue "C6" "5,022"

# New row 7
Set-TextValue "A7" "2024-10-25 16:48"
Set-TextValue "B7" "2"
Set-TextValue "C7" "18,832"
Set-TextValue "D7" "79,648"

# New row 8
Set-TextValue "A8" "2024-10-10 17:23"
Set-TextValue "B8" "2"
Set-TextValue "C8" "21,971"
Set-TextValue "D8" "79,648"

# New row 9
Set-TextValue "A9" "2024-09-30 17:45"
Set-TextValue "B9" "2"
Set-TextValue "C9" "22,599"
Set-TextValue "D9" "79,648"
